$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: move the existing "Customer Count" values (column B) into the new column E ---
$ws.Range("E2").Value2 = 279
$ws.Range("E3").Value2 = 103
$ws.Range("E4").Value2 = 68
$ws.Range("E5").Value2 = 343

# --- Step 2: overwrite column B with the new Recency_Mean values ---
$ws.Range("B2").Value2 = 71.14
$ws.Range("B3").Value2 = 547.89
$ws.Range("B4").Value2 = 120.71
$ws.Range("B5").Value2 = 98.81999999999999

# --- Step 3: rename headers ---
$ws.Range("A1").Value2 = "Cluster"
$ws.Range("B1").Value2 = "Recency_Mean"
$ws.Range("C1").Value2 = "Frequency_Mean"
$ws.Range("D1").Value2 = "Monetary_Mean"
$ws.Range("E1").Value2 = "Customer_Count"

# --- Step 4: strip the bold/border/center formatting from the header row so it uses the default style ---
$ws.Range("A1:E1").ClearFormats()

# --- Step 5: make sure C/D values still match exactly (they are unchanged, but restate to be safe) ---
$ws.Range("C2").Value2 = 8.5
$ws.Range("D2").Value2 = 3226.84
$ws.Range("C3").Value2 = 3.71
$ws.Range("D3").Value2 = 1432.86
$ws.Range("C4").Value2 = 8.210000000000001
$ws.Range("D4").Value2 = 9236.450000000001
$ws.Range("C5").Value2 = 4.7
$ws.Range("D5").Value2 = 1707.24
